$wb = $excel.ActiveWorkbook

# Hoja2 (Sheet2): change A1 text and add a new B1 cell, then make this the
# active/selected sheet (as in the final workbook).
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("B1").Value = "call"
$ws2.Range("A1").Value = "gateway"
$ws2.Activate()
